$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column A entirely (old TAXON-count column that duplicated the GENE
# count in column F) - this shifts B:F left to A:E.
$ws.Range("A1").EntireColumn.Delete()

# Rename the header that used to read "MODEL_CONDITION" (now in column D
# after the shift) to "MODELCONDITION".
$ws.Range("D1").Value = "MODELCONDITION"
